# Insert a new record row at row 256 (pushes existing rows 256-366 down to 257-367)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("256:256").Insert()

$ws.Cells.Item(256, 1).Value = 10
$ws.Cells.Item(256, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(256, 3).Value = "La Araucanía"
$ws.Cells.Item(256, 4).Value = 44875
$ws.Cells.Item(256, 5).Value = 9
$ws.Cells.Item(256, 6).Value = 100114013
$ws.Cells.Item(256, 7).Value = "Zanahoria"
$ws.Cells.Item(256, 8).Value = "Sin especificar"
$ws.Cells.Item(256, 9).Value = "Primera"
$ws.Cells.Item(256, 10).Value = 380
$ws.Cells.Item(256, 11).Value = 14000
$ws.Cells.Item(256, 12).Value = 14000
$ws.Cells.Item(256, 13).Value = 14000
$ws.Cells.Item(256, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(256, 15).Value = "Región del Bíobío"
$ws.Cells.Item(256, 16).Value = 700
$ws.Cells.Item(256, 17).Value = 20
$ws.Cells.Item(256, 18).Value = "Hortaliza"
